# Fruta / hortaliza, semanal
# Insert a new weekly record as row 193 on the "Durazno" sheet,
# shifting the existing rows 193:294 down to 194:295.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 193; rows below (193-294) shift to 194-295
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new record's data
$ws.Cells.Item(193, 1).Value  = 10
$ws.Cells.Item(193, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value  = "La Araucanía"
$ws.Cells.Item(193, 4).Value  = 44875
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 5).Value  = 9
$ws.Cells.Item(193, 6).Value  = "Fruta"
$ws.Cells.Item(193, 7).Value  = 100103
$ws.Cells.Item(193, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(193, 9).Value  = 100103004
$ws.Cells.Item(193, 10).Value = "Durazno"
$ws.Cells.Item(193, 11).Value = "Early Majestic"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 650
$ws.Cells.Item(193, 14).Value = 20000
$ws.Cells.Item(193, 15).Value = 22000
$ws.Cells.Item(193, 16).Value = 21077
$ws.Cells.Item(193, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(193, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(193, 19).Value = 2108
$ws.Cells.Item(193, 20).Value = 10
